# Lessons-learned template: turn the single empty paragraph (which only
# carries the _GoBack bookmark) into six bold, Times New Roman, double
# spaced prompt paragraphs. The first paragraph keeps the bookmark.

$d = $word.ActiveDocument
$p1 = $d.Paragraphs.Item(1)

function RunPropsXml([bool]$eastAsia, [string]$color) {
    if ($eastAsia) {
        $fonts = '<w:rFonts w:ascii="Times New Roman" w:eastAsia="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>'
    } else {
        $fonts = '<w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>'
    }
    $colorXml = ''
    if ($color) {
        $colorXml = '<w:color w:val="' + $color + '"/>'
    }
    return '<w:rPr>' + $fonts + '<w:b/>' + $colorXml + '<w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr>'
}

function ParaXml([string]$text, [bool]$eastAsia, [string]$color, [bool]$withBookmark) {
    $rpr = RunPropsXml $eastAsia $color
    $xml = '<w:p><w:pPr><w:spacing w:line="480" w:lineRule="auto"/>' + $rpr + '</w:pPr>'
    $xml += '<w:r>' + $rpr + '<w:t>' + $text + '</w:t></w:r>'
    if ($withBookmark) {
        $xml += '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>'
    }
    $xml += '</w:p>'
    return $xml
}

$paragraphs = @(
    @{ Text = "What went well?"; EastAsia = $false; Color = $null },
    @{ Text = "What went poorly?"; EastAsia = $false; Color = $null },
    @{ Text = "What would you do differently?"; EastAsia = $false; Color = $null },
    @{ Text = "If this opportunity presented itself again, would you do it?"; EastAsia = $true; Color = "111111" },
    @{ Text = "Were you profitable?"; EastAsia = $false; Color = $null },
    @{ Text = "Goals for next time?"; EastAsia = $false; Color = $null }
)

$bodyXml = ""
for ($i = 0; $i -lt $paragraphs.Count; $i++) {
    $item = $paragraphs[$i]
    $isFirst = ($i -eq 0)
    $bodyXml += ParaXml $item.Text $item.EastAsia $item.Color $isFirst
}

$packageXml = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $bodyXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# Replacing paragraph 1's range (which today is just the bookmarked,
# empty paragraph) with the whole six-paragraph block in one shot avoids
# leaving a stray empty paragraph behind, and keeps the _GoBack bookmark
# anchored in the first ("What went well?") paragraph as in the original.
$p1.Range.InsertXML($packageXml)
